$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.762.87'
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").Value = '1.884.63'
$ws.Range("E3").Value = '  -5.00%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.49'
$ws.Range("E5").Value = '  -1.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3794'
$ws.Range("E8").Value = '  -3.78%  '
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07704'
$ws.Range("E10").Value = '  -2.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9602'
$ws.Range("E11").Value = '  -4.10%  '
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("D13").Value = '1.889.61'
$ws.Range("E13").Value = '  -4.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.938'
$ws.Range("E14").Value = '  -3.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.653'
$ws.Range("E15").Value = '  -3.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06986'
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '83.07'
$ws.Range("E18").Value = '  -6.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009470'
$ws.Range("E19").Value = '  -4.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.57'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = '28.713.91'
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.315'
$ws.Range("E23").Value = '  -4.27%  '
$ws.Range("E24").Value = '  -3.47%  '
$ws.Range("D25").Value = '2.123.32'
$ws.Range("E25").Value = '  -4.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.084'
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.16'
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.93'
$ws.Range("E28").Value = '  -3.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.585'
$ws.Range("E29").Value = '  -6.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.96'
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.798'
$ws.Range("E31").Value = '  -6.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09237'
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8419'
$ws.Range("E33").Value = '  -5.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.059'
$ws.Range("E34").Value = '  -4.11%  '
$ws.Range("E35").Value = '  -8.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.983'
$ws.Range("E36").Value = '  -5.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05653'
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.138'
$ws.Range("E38").Value = '  -3.35%  '
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02024'
$ws.Range("E40").Value = '  -4.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.406'
$ws.Range("E41").Value = '  -6.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5469'
$ws.Range("E42").Value = '  -5.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1744'
$ws.Range("E43").Value = '  -4.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000002974'
$ws.Range("E44").Value = '  -30.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.111'
$ws.Range("E45").Value = '  -7.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.690'
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5139'
$ws.Range("E47").Value = '  -4.47%  '
$ws.Range("E48").Value = '  -6.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06791'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.058'
$ws.Range("E50").Value = '  -5.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.44'
$ws.Range("E51").Value = '  -2.62%  '
